$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string text referenced by E3/L3 ("ap_clk_IBUF" -> "ap_clk_IBUF_BUFG")
$ws.Range("E3").Value = "ap_clk_IBUF_BUFG"
$ws.Range("L3").Value = "ap_clk_IBUF_BUFG"

# Update numeric values
$ws.Range("A2").Value = 0.0012505515478551388
$ws.Range("A3").Value = 0.0012505515478551388
$ws.Range("D3").Value = 9.324859619140625
$ws.Range("G3").Value = 19.671045303344727
$ws.Range("K3").Value = 10.667463302612305
$ws.Range("N3").Value = 39.26985168457031
$ws.Range("R3").Value = 8.95522403717041
